$wb = $excel.ActiveWorkbook

# hunk 0: ALC!row33
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 559.75
$ws.Range("I33").Value = 566.5833
$ws.Range("K33").Value = 566.5833
$ws.Range("M33").Value = -337.5833

# hunk 1: ALC!row49
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H49").Value = 399.5
$ws.Range("I49").Value = 300
$ws.Range("K49").Value = 900
$ws.Range("M49").Value = -764

# hunk 2: ALC!row64
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3790
$ws.Range("J64").Value = 3980
$ws.Range("L64").Value = 3980
$ws.Range("N64").Value = -4476

# hunk 3: ALC!row67
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 3790
$ws.Range("J67").Value = 3980
$ws.Range("L67").Value = 3980
$ws.Range("N67").Value = -5696

# hunk 4: ALC!row112
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 2294.65
$ws.Range("I112").Value = 833.3333
$ws.Range("J112").Value = 2552.5293
$ws.Range("K112").Value = 2499.9999
$ws.Range("L112").Value = 7657.5879
$ws.Range("M112").Value = -1391.9999
$ws.Range("N112").Value = -9873.5879

# hunk 5: ALC!row129
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 820.3684
$ws.Range("I129").Value = 542
$ws.Range("J129").Value = 872.5625
$ws.Range("K129").Value = 1626
$ws.Range("L129").Value = 2617.6875
$ws.Range("M129").Value = 3374
$ws.Range("N129").Value = -12617.6875

# hunk 6: ALC!row135
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 55557270
$ws.Range("I135").Value = 526.2
$ws.Range("J135").Value = 333341000
$ws.Range("K135").Value = 4735.8
$ws.Range("L135").Value = 3000069000
$ws.Range("M135").Value = -2200.8
$ws.Range("N135").Value = -3000074070

# hunk 7: ALC!row137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1635.875
$ws.Range("I137").Value = 1194.2142
$ws.Range("J137").Value = 1979.3889
$ws.Range("K137").Value = 3582.6426
$ws.Range("L137").Value = 5938.1667
$ws.Range("M137").Value = -1032.6426
$ws.Range("N137").Value = -11038.1667

# hunk 8: ALC!row138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 459613.88
$ws.Range("I138").Value = 1444.2
$ws.Range("J138").Value = 545520.7
$ws.Range("K138").Value = 4332.6
$ws.Range("L138").Value = 1636562.1
$ws.Range("M138").Value = 807.3999999999996
$ws.Range("N138").Value = -1646842.1

# hunk 9: ALC!row141
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 4517
$ws.Range("I141").Value = 4517
$ws.Range("K141").Value = 13551
$ws.Range("M141").Value = -8371

# hunk 10: ARM!row63
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 21278404
$ws.Range("I63").Value = 1640.5625
$ws.Range("J63").Value = 66668830
$ws.Range("K63").Value = 1640.5625
$ws.Range("L63").Value = 66668830
$ws.Range("M63").Value = -954.5625
$ws.Range("N63").Value = -66670202

# hunk 11: ARM!row66
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 21278404
$ws.Range("I66").Value = 1640.5625
$ws.Range("J66").Value = 66668830
$ws.Range("K66").Value = 8202.8125
$ws.Range("L66").Value = 333344150
$ws.Range("M66").Value = -4770.8125
$ws.Range("N66").Value = -333351014

# hunk 12: ARM!row74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1992.4706
$ws.Range("I74").Value = 1124.7273
$ws.Range("J74").Value = 3583.3333
$ws.Range("K74").Value = 1124.7273
$ws.Range("L74").Value = 3583.3333
$ws.Range("M74").Value = -250.7273
$ws.Range("N74").Value = -5331.3333

# hunk 13: ARM!row77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1992.4706
$ws.Range("I77").Value = 1124.7273
$ws.Range("J77").Value = 3583.3333
$ws.Range("K77").Value = 5623.636500000001
$ws.Range("L77").Value = 17916.6665
$ws.Range("M77").Value = -1255.636500000001
$ws.Range("N77").Value = -26652.6665

# hunk 14: ARM!row118
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").ClearContents()

# hunk 15: ARM!row122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2255.4666
$ws.Range("I122").Value = 2083
$ws.Range("J122").Value = 2452.5715
$ws.Range("K122").Value = 6249
$ws.Range("L122").Value = 7357.7145
$ws.Range("M122").Value = -3799
$ws.Range("N122").Value = -12257.7145

# hunk 16: BSM!row134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1162
$ws.Range("I134").Value = 939.0833
$ws.Range("J134").Value = 2499.5
$ws.Range("K134").Value = 2817.2499
$ws.Range("L134").Value = 7498.5
$ws.Range("M134").Value = -282.2498999999998
$ws.Range("N134").Value = -12568.5

# hunk 17: CRP!row31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1440.0238
$ws.Range("I31").Value = 1322.6154
$ws.Range("J31").Value = 2966.3333
$ws.Range("K31").Value = 1322.6154
$ws.Range("L31").Value = 2966.3333
$ws.Range("M31").Value = -1027.6154
$ws.Range("N31").Value = -3556.3333

# hunk 18: CRP!row34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1440.0238
$ws.Range("I34").Value = 1322.6154
$ws.Range("J34").Value = 2966.3333
$ws.Range("K34").Value = 1322.6154
$ws.Range("L34").Value = 2966.3333
$ws.Range("M34").Value = -1120.6154
$ws.Range("N34").Value = -3370.3333

# hunk 19: CRP!row122
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1850.5
$ws.Range("I122").Value = 1800.6666
$ws.Range("K122").Value = 5401.9998
$ws.Range("M122").Value = -2951.9998

# hunk 20: CUL!row40
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 379.4
$ws.Range("I40").Value = 148.2
$ws.Range("J40").Value = 495
$ws.Range("K40").Value = 592.8
$ws.Range("L40").Value = 1980
$ws.Range("M40").Value = -523.8
$ws.Range("N40").Value = -2118

# hunk 21: CUL!row62
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 100
$ws.Range("I62").Value = 100
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 300
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("M62").Value = 386

# hunk 22: CUL!row65
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H65").Value = 100
$ws.Range("I65").Value = 100
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 900
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("M65").Value = 2532

# hunk 23: CUL!row113
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 557.8409
$ws.Range("I113").Value = 474.125
$ws.Range("J113").Value = 605.6786
$ws.Range("K113").Value = 1422.375
$ws.Range("L113").Value = 1817.0358
$ws.Range("M113").Value = 747.625
$ws.Range("N113").Value = -6157.0358

# hunk 24: CUL!row131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 16950092
$ws.Range("I131").Value = 125000424
$ws.Range("J131").Value = 1020.7647
$ws.Range("K131").Value = 375001272
$ws.Range("L131").Value = 3062.2941
$ws.Range("M131").Value = -374996232
$ws.Range("N131").Value = -13142.2941

# hunk 25: GSM!row51
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H51").Value = 29999
$ws.Range("J51").Value = 29999
$ws.Range("L51").Value = 29999
$ws.Range("N51").Value = -31017

# hunk 26: LTW!row46
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 4709.9
$ws.Range("J46").Value = 5155.4443
$ws.Range("L46").Value = 5155.4443
$ws.Range("N46").Value = -5531.4443

# hunk 27: LTW!row100
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 1300
$ws.Range("I100").Value = 1200
$ws.Range("J100").Value = 1360
$ws.Range("K100").Value = 1200
$ws.Range("L100").Value = 1360
$ws.Range("M100").Value = -659
$ws.Range("N100").Value = -2442

# hunk 28: LTW!row122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 35716212
$ws.Range("I122").Value = 62501868
$ws.Range("J122").Value = 2001.3334
$ws.Range("K122").Value = 187505604
$ws.Range("L122").Value = 6004.0002
$ws.Range("M122").Value = -187503154
$ws.Range("N122").Value = -10904.0002

# hunk 29: WVR!row8
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H8").Value = 85004
$ws.Range("J8").Value = 85004
$ws.Range("L8").Value = 85004
$ws.Range("N8").Value = -85284

# hunk 30: WVR!row11
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H11").Value = 800
$ws.Range("J11").Value = 800
$ws.Range("L11").Value = 800
$ws.Range("N11").Value = -1084

# hunk 31: WVR!row119
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H119").Value = 15000
$ws.Range("J119").Value = 15000
$ws.Range("L119").Value = 15000
$ws.Range("N119").Value = -24676
